$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 687-688, shifting the existing data (old rows
# 687..781) down to 689..783.
$ws.Rows("687:688").Insert()

# The "constant" attribute columns (A..K) are identical for every row in
# this Pina/Caramelo block, so copy them from the row that used to be 687
# (now sitting at 689) into the two freshly inserted rows.
for ($c = 1; $c -le 11; $c++) {
    $v = $ws.Cells.Item(689, $c).Value2
    $ws.Cells.Item(687, $c).Value2 = $v
    $ws.Cells.Item(688, $c).Value2 = $v
}

# New row 687: "Primera" quality quote.
$ws.Cells.Item(687, 4).Value2 = 45142
$ws.Cells.Item(687, 12).Value2 = "Primera"
$ws.Cells.Item(687, 13).Value2 = 95
$ws.Cells.Item(687, 14).Value2 = 24000
$ws.Cells.Item(687, 15).Value2 = 24000
$ws.Cells.Item(687, 16).Value2 = 24000
$ws.Cells.Item(687, 17).Value2 = "`$/caja 12 unidades"
$ws.Cells.Item(687, 18).Value2 = "Ecuador"
$ws.Cells.Item(687, 19).Value2 = 2000
$ws.Cells.Item(687, 20).Value2 = 12

# New row 688: "Segunda" quality quote.
$ws.Cells.Item(688, 4).Value2 = 45142
$ws.Cells.Item(688, 12).Value2 = "Segunda"
$ws.Cells.Item(688, 13).Value2 = 65
$ws.Cells.Item(688, 14).Value2 = 24000
$ws.Cells.Item(688, 15).Value2 = 24000
$ws.Cells.Item(688, 16).Value2 = 24000
$ws.Cells.Item(688, 17).Value2 = "`$/caja 14 unidades"
$ws.Cells.Item(688, 18).Value2 = "Ecuador"
$ws.Cells.Item(688, 19).Value2 = 1714
$ws.Cells.Item(688, 20).Value2 = 14

# Match the date-format style used by the rest of column D.
$ws.Range("D687:D688").NumberFormat = $ws.Range("D689").NumberFormat
